$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.879.46'
$ws.Range('E2').Value = '  +2.64%  '
$ws.Range('D3').Value = '3.033.04'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.71%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.027.86'
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.151'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.463'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.72%  '
$ws.Range('E13').Value = '  +2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.20%  '
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '3.537.18'
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.08'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D18').Value = '62.899.94'
$ws.Range('E18').Value = '  +2.70%  '
$ws.Range('D19').Value = '3.035.62'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.697'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.17%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +4.59%  '
$ws.Range('E30').Value = '  +10.09%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.59'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('D35').Value = '0.0₃0858'
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('E36').Value = '  +2.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.130'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('E43').Value = '  +14.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.87'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '391.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0361'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.87%  '
$ws.Range('D47').Value = '2.721.87'
$ws.Range('E47').Value = '  +1.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.71%  '
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.72%  '
